$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new price-record row is inserted at row 545, pushing all existing
# records (545..651) down by one row (to 546..652).
$ws.Rows("545:545").Insert()

$newRow = 545
$ws.Cells.Item($newRow, 1).Value = 7
$ws.Cells.Item($newRow, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($newRow, 3).Value = "Ñuble"
$ws.Cells.Item($newRow, 4).Value = 44995
$ws.Cells.Item($newRow, 5).Value = 16
$ws.Cells.Item($newRow, 6).Value = 100112020
$ws.Cells.Item($newRow, 7).Value = "Tomate"
$ws.Cells.Item($newRow, 8).Value = "Semiduro"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 600
$ws.Cells.Item($newRow, 11).Value = 7500
$ws.Cells.Item($newRow, 12).Value = 8000
$ws.Cells.Item($newRow, 13).Value = 7750
$ws.Cells.Item($newRow, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item($newRow, 15).Value = "Región del Maule"
$ws.Cells.Item($newRow, 16).Value = 431
$ws.Cells.Item($newRow, 17).Value = 18
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"

# Match the date-format style used by the other rows' "Fecha" column.
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow + 1, 4).NumberFormat
